$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (ROC1), Outro2: update "tree" -> "TREE" ---
$ws.Range('G3').Value = 'You find a group of the {character:baron:baron}''s men in the tavern drinking. As one of them turns to you, your first blow lands squarely on his jaw. An all-out brawl breaks out. You give as good as you get, but in the end there are too many of them.

"I''ll take care of {objPronoun}," you hear a voice say. Rough hands carry you out of town into the nearby {location:nearby:name:pathtobaron}. You pass out.{|GOTO:pathtobaron|}

You wake up the next day, deep in {location:current:namewiththe}. A {character:mentor:sexAge:ranger} sits under a nearby TREE. As you open your eyes, {character:ranger:subPronoun} says, "Hope you appreciate I saved your life. That was kind of a stupid fight to start, you know?"'

# --- Row 4 (MTM1): change condition separator | -> & ---
$ws.Range('B4').Value = 'character:baron&character:ranger'

# --- Row 6 (BOTW1): change condition separator and add location condition ---
$ws.Range('B6').Value = 'character:baron&character:ranger&location:current:forest'

# --- Row 7 (BOTW1a): clear condition entirely ---
$ws.Range('B7').ClearContents()

# --- Row 8 (ROT1): append SET flag to message, and update scamper text ---
$ws.Range('C8').Value = 'The next day, a small creature runs across your path. You almost step on it.

You squat down to look at it. It seems like some sort of squirrel, but it''s eyes are ruby red. It sits still, staring at you, alert.{|SET:squirrel:exists|}'
$ws.Range('F8').Value = 'The squirrel seems appreciative, and chirps happily. It looks deep into your eyes, then it scampers off into {location:current:namewiththe}. '

# --- Row 10 (ROT3): add new condition ---
$ws.Range('B10').Value = 'location:current:forest&squirrel:exists'

# --- Row 11 (MWG1): change condition separator and add location condition ---
$ws.Range('B11').Value = 'character:baron&character:ranger&location:current:forest'

# --- Row 13 (AWF1): change condition separator ---
$ws.Range('B13').Value = 'character:baron&item:map'

# --- Row 14 (AWF2): change condition separator ---
$ws.Range('B14').Value = 'character:baron&item:map&character:stolen'

# --- Row 15 (AWF3): change condition separator ---
$ws.Range('B15').Value = 'character:baron&character:stolen'

# --- Row 17 (AWF1a): clear condition entirely ---
$ws.Range('B17').ClearContents()

# --- Row 18 (A1): change condition separator ---
$ws.Range('B18').Value = 'character:stolen&character:baron'

# --- Update selected/active cell and scroll position to match authored view ---
$ws.Range('G3').Select()
